# Auto-generated Word COM-interop script implementing the content rotation described by the diff.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $ok = $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $find"
    }
}

# ---- Phase 1: stash each moved text block behind a unique placeholder ----
$v1 = "Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental."
$v2 = "PLACEHOLDER_A1"
Replace-Text $v1 $v2   # P6
$v3 = "9146830 - Danúbia Caporusso Bargos"
$v4 = "PLACEHOLDER_A2"
Replace-Text $v3 $v4   # P9run1
$v5 = "Provas e/ou exercícios dirigidos"
$v6 = "PLACEHOLDER_A3"
Replace-Text $v5 $v6   # P17run6
$v7 = "Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas"
$v8 = "PLACEHOLDER_A4"
Replace-Text $v7 $v8   # P17run2
$v9 = "Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental."
$v10 = "PLACEHOLDER_A5"
Replace-Text $v9 $v10   # P11
$v11 = "Provide students with knowledge of the fundamentals of Environmental Education using as basis the current environmental problems. To develop practical activities integrated to the region. Guide the development of projects related to Environmental Education and Management"
$v12 = "PLACEHOLDER_B1"
Replace-Text $v11 $v12   # P7
$v13 = "General considerations on environmental problem. Evolution of environmental questions in Brazil and in the world. Education and Environmental Management. Development and monitoring of environmental education projects."
$v14 = "PLACEHOLDER_B2"
Replace-Text $v13 $v14   # P12
$v15 = "Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental."
$v16 = "PLACEHOLDER_C1"
Replace-Text $v15 $v16   # P14
$v17 = "5817650 - Érica Leonor Romão"
$v18 = "PLACEHOLDER_C2"
Replace-Text $v17 $v18   # P9run2
$v19 = "CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006." + [char]11 + "CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1." + [char]11 + "DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p." + [char]11 + "DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p." + [char]11 + "DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000." + [char]11 + "GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006." + [char]11 + "JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p." + [char]11 + "JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p " + [char]11 + "LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165p" + [char]11 + "PHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3)."
$v20 = "PLACEHOLDER_C3"
Replace-Text $v19 $v20   # P19
$v21 = "Nota Final: NF ≥ 5,0"
$v22 = "PLACEHOLDER_C4"
Replace-Text $v21 $v22   # P17run4

# ---- Phase 2: drop each text block into its final destination ----
$v23 = "PLACEHOLDER_A1"
$v24 = "Considerações gerais sobre a problemática ambiental. Evolução das questões ambientais no Brasil e no mundo. Educação e Gestão Ambiental. Elaboração e acompanhamento de projetos de educação ambiental."
Replace-Text $v23 $v24
$v25 = "PLACEHOLDER_A2"
$v26 = "Propiciar ao discente conhecimento dos fundamentos da Educação Ambiental utilizando como base os problemas ambientais da atualidade. Desenvolver atividades práticas integradas à região. Orientar o desenvolvimento de projetos relacionados à Gestão e Educação Ambiental."
Replace-Text $v25 $v26
$v27 = "PLACEHOLDER_A3"
$v28 = "9146830 - Danúbia Caporusso Bargos"
Replace-Text $v27 $v28
$v29 = "PLACEHOLDER_A4"
$v30 = "Provas e/ou exercícios dirigidos"
Replace-Text $v29 $v30
$v31 = "PLACEHOLDER_A5"
$v32 = "Avaliação baseada em provas, exercícios, projetos, seminários e outras formas de avaliação, sendo a nota final correspondente a média ponderada das notas atribuídas às avaliações aplicadas"
Replace-Text $v31 $v32
$v33 = "PLACEHOLDER_B1"
$v34 = "General considerations on environmental problem. Evolution of environmental questions in Brazil and in the world. Education and Environmental Management. Development and monitoring of environmental education projects."
Replace-Text $v33 $v34
$v35 = "PLACEHOLDER_B2"
$v36 = "Provide students with knowledge of the fundamentals of Environmental Education using as basis the current environmental problems. To develop practical activities integrated to the region. Guide the development of projects related to Environmental Education and Management"
Replace-Text $v35 $v36
$v37 = "PLACEHOLDER_C1"
$v38 = "Nota Final: NF ≥ 5,0"
Replace-Text $v37 $v38
$v39 = "PLACEHOLDER_C2"
$v40 = "Sociedade, natureza e desenvolvimento. A relação degradação ambiental-qualidade de vida. Meio ambiente e cidadania. Percepção e Interpretação ambiental. Meio ambiente e representação social. Histórico da educação ambiental e conceitos de meio ambiente; Conceitos, princípios e pensamentos norteadores da Educação Ambiental. A questão ambiental e as conferências mundiais de meio ambiente. O movimento ambientalista e o histórico da EA no Brasil e no mundo; A Agenda 21 e educação ambiental. A política nacional de educação ambiental (pnea) e legislação correlata: A abordagem interdisciplinar da educação ambiental; Educação como instrumento de Gestão Ambiental. Educação ambiental nas empresas e o Sistema de Gestão Ambiental. Projetos, reflexões e práticas da Educação Ambiental. Análise e vivências de experiências práticas de educação ambiental em diferentes contextos. Metodologia de projetos, oficinas e capacitação em educação ambiental."
Replace-Text $v39 $v40
$v41 = "PLACEHOLDER_C3"
$v42 = "5817650 - Érica Leonor Romão"
Replace-Text $v41 $v42
$v43 = "PLACEHOLDER_C4"
$v44 = "CARVALHO, I. C. M.; Educação Ambiental e formação do sujeito ecológico. São Paulo: Cortez, 2006." + [char]11 + "CINQUETTI, H. C. S.; LOGAREZZI, A. (Org.). Consumo e Resíduo - Fundamentos para o trabalho educativo. 1 ed. São Carlos: EdUFSCar, 2006, v. 1." + [char]11 + "DIAS, G. F. Dinâmica e instrumentação para educação ambiental. 1. ed. São Paulo: Gaia, 2010. v. 1. 216p." + [char]11 + "DIAS, G. F. Educação e Gestão Ambiental. 1. ed. São Paulo: Editora Gaia Ltda, 2006. v. 1. 118p." + [char]11 + "DIAS, G. F. Educação Ambiental: princípios e práticas. 6a ed. São Paulo: Gaia, 2000." + [char]11 + "GUIMARÃES, M. (org.) Caminhos da educação ambiental: da forma à ação. Campinas, SP: Papirus, 2006." + [char]11 + "JACOBI, Pedro Roberto, MONTEIRO,F. M ; FERNANDES, M. L. B. . Educação e Sustentabilidade- caminhos e práticas para uma educação transformadora. São Paulo: Evoluir Cultural, 2009. v. 01. 108p." + [char]11 + "JACOBI, Pedro Roberto OLIVEIRA, F. C. J. F. (Org.). Educação, Meio Ambiente e Cidadania - reflexões e experiências. São Paulo: SMA/CEAM, 1998. 121p " + [char]11 + "LOUREIRO, C. F. B. Trajetória e fundamentos da educação ambiental. 4. ed. São Paulo: Cortez editora, 2012. 165p" + [char]11 + "PHILIPPI JR., A & PELICIONI, M. C. F. (Eds). 2005. Educação ambiental e sustentabilidade. Barueri SP: Manole. 878p. (Coleção Ambiental, 3)."
Replace-Text $v43 $v44

Write-Output "done"